# Apply "Fixed naive component forecaster bug - Presentation state 11.02" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear C2 (no longer has a forecast value) and update E2 precision
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 19.1981274365808

# Row 3
$ws.Range("E3").Value = -9.964084247724713

# Row 4
$ws.Range("C4").Value = -5.440152375872276
$ws.Range("E4").Value = -14.43639438706736

# Row 5
$ws.Range("C5").Value = 9.349082908138474
$ws.Range("E5").Value = 27.15801420548431

# Row 6
$ws.Range("C6").Value = 0.5389546843749926

# Row 7
$ws.Range("C7").Value = -4.232836797447703
$ws.Range("E7").Value = -8.093075920532211

# Row 8
$ws.Range("E8").Value = 15.37760125310901

# Row 9
$ws.Range("C9").Value = 1.913895196850146
$ws.Range("E9").Value = 6.974907992137958

# Row 10
$ws.Range("C10").Value = 4.861901970954019

# Row 11
$ws.Range("C11").Value = 4.115488239647735
$ws.Range("E11").Value = 9.52305004616103

# Row 12
$ws.Range("C12").Value = 4.073811422566442
$ws.Range("E12").Value = 8.549566886636839

# Row 13
$ws.Range("E13").Value = 7.819356632099983

# Row 15
$ws.Range("C15").Value = 3.285232806602423
$ws.Range("E15").Value = 4.648946574958668

# Row 17
$ws.Range("C17").Value = 1.16693824877212
$ws.Range("E17").Value = 16.68718678695833

# Row 18
$ws.Range("C18").Value = 1.758584501904181

# Row 19
$ws.Range("C19").Value = 0.5908161348962437
